$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2303.0667
$ws.Range("J19").Value = 1564.4445
$ws.Range("L19").Value = 1564.4445
$ws.Range("N19").Value = -1914.4445

$ws.Range("H58").Value = 1387.5
$ws.Range("J58").Value = 1625
$ws.Range("L58").Value = 4875
$ws.Range("N58").Value = -5175

$ws.Range("H64").Value = 4179.75
$ws.Range("J64").Value = 5850
$ws.Range("L64").Value = 5850
$ws.Range("N64").Value = -6346

$ws.Range("H67").Value = 4179.75
$ws.Range("J67").Value = 5850
$ws.Range("L67").Value = 5850
$ws.Range("N67").Value = -7566

$ws.Range("H76").Value = 3492.6667
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 3492.6667
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()

$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()

$ws.Range("H97").Value = 1804
$ws.Range("J97").Value = 1804
$ws.Range("L97").Value = 5412
$ws.Range("N97").Value = -6404

$ws.Range("H100").Value = 2791.7693
$ws.Range("I100").Value = 2523.125
$ws.Range("J100").Value = 3221.6
$ws.Range("K100").Value = 2523.125
$ws.Range("L100").Value = 3221.6
$ws.Range("M100").Value = -1982.125
$ws.Range("N100").Value = -4303.6

$ws.Range("H107").Value = 1064.8
$ws.Range("I107").Value = 1064.8
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1064.8
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 855.2
$ws.Range("N107").ClearContents()

$ws.Range("H135").Value = 626.1613
$ws.Range("I135").Value = 626.1613
$ws.Range("K135").Value = 5635.4517
$ws.Range("M135").Value = -3100.4517

$ws.Range("H138").Value = 1875.8485
$ws.Range("I138").Value = 1162.7916
$ws.Range("J138").Value = 3777.3333
$ws.Range("K138").Value = 3488.3748
$ws.Range("L138").Value = 11331.9999
$ws.Range("M138").Value = 1651.6252
$ws.Range("N138").Value = -21611.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 20000
$ws.Range("J23").Value = 20000
$ws.Range("L23").Value = 20000
$ws.Range("N23").Value = -20518

$ws.Range("H32").Value = 34392.81
$ws.Range("I32").Value = 19645.852
$ws.Range("K32").Value = 19645.852
$ws.Range("M32").Value = -19358.852

$ws.Range("H61").Value = 1975.5
$ws.Range("I61").Value = 1554.8334
$ws.Range("K61").Value = 1554.8334
$ws.Range("M61").Value = -1342.8334

$ws.Range("H110").Value = 1958.7778
$ws.Range("I110").Value = 2057.9333
$ws.Range("J110").Value = 1463
$ws.Range("K110").Value = 2057.9333
$ws.Range("L110").Value = 1463
$ws.Range("M110").Value = -12.93330000000014
$ws.Range("N110").Value = -5553

$ws.Range("H132").Value = 16608
$ws.Range("I132").Value = 20466.363
$ws.Range("J132").Value = 5997.5
$ws.Range("K132").Value = 61399.08900000001
$ws.Range("L132").Value = 17992.5
$ws.Range("M132").Value = -58869.08900000001
$ws.Range("N132").Value = -23052.5

$ws.Range("H136").Value = 1975.5
$ws.Range("I136").Value = 1554.8334
$ws.Range("K136").Value = 4664.5002
$ws.Range("M136").Value = -2114.5002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 8902.416999999999
$ws.Range("I107").Value = 2768.8
$ws.Range("J107").Value = 39570.5
$ws.Range("K107").Value = 2768.8
$ws.Range("L107").Value = 39570.5
$ws.Range("M107").Value = -848.8000000000002
$ws.Range("N107").Value = -43410.5

$ws.Range("H134").Value = 965.14813
$ws.Range("I134").Value = 938.4
$ws.Range("K134").Value = 2815.2
$ws.Range("M134").Value = -280.1999999999998

$ws.Range("H138").Value = 50198.8
$ws.Range("J138").Value = 50198.8
$ws.Range("L138").Value = 50198.8
$ws.Range("N138").Value = -60478.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 380
$ws.Range("I22").Value = 380
$ws.Range("K22").Value = 380
$ws.Range("M22").Value = -30

$ws.Range("H31").Value = 1840.7727
$ws.Range("I31").Value = 1409.303
$ws.Range("K31").Value = 1409.303
$ws.Range("M31").Value = -1114.303

$ws.Range("H34").Value = 1840.7727
$ws.Range("I34").Value = 1409.303
$ws.Range("K34").Value = 1409.303
$ws.Range("M34").Value = -1207.303

$ws.Range("H92").Value = 30601
$ws.Range("J92").Value = 30601
$ws.Range("L92").Value = 30601
$ws.Range("N92").Value = -35593

$ws.Range("H132").Value = 3390.8333
$ws.Range("I132").Value = 3203.2307
$ws.Range("K132").Value = 9609.6921
$ws.Range("M132").Value = -7079.6921

$ws.Range("H134").Value = 3756.923
$ws.Range("I134").Value = 3531.818
$ws.Range("J134").Value = 4995
$ws.Range("K134").Value = 10595.454
$ws.Range("L134").Value = 14985
$ws.Range("M134").Value = -8060.454000000002
$ws.Range("N134").Value = -20055

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 6652.625
$ws.Range("I69").Value = 9333.333000000001
$ws.Range("J69").Value = 5044.2
$ws.Range("K69").Value = 27999.999
$ws.Range("L69").Value = 15132.6
$ws.Range("M69").Value = -27188.999
$ws.Range("N69").Value = -16754.6

$ws.Range("H72").Value = 6652.625
$ws.Range("I72").Value = 9333.333000000001
$ws.Range("J72").Value = 5044.2
$ws.Range("K72").Value = 83999.997
$ws.Range("L72").Value = 45397.8
$ws.Range("M72").Value = -79943.997
$ws.Range("N72").Value = -53509.8

$ws.Range("H107").Value = 853.4286
$ws.Range("I107").Value = 427.33334
$ws.Range("J107").Value = 1055.2632
$ws.Range("K107").Value = 1282.00002
$ws.Range("L107").Value = 3165.7896
$ws.Range("M107").Value = 637.9999800000001
$ws.Range("N107").Value = -7005.7896

$ws.Range("H128").Value = 318344
$ws.Range("I128").Value = 318344
$ws.Range("K128").Value = 955032
$ws.Range("M128").Value = -950052

$ws.Range("H134").Value = 1207.6
$ws.Range("I134").Value = 1207.6
$ws.Range("K134").Value = 3622.8
$ws.Range("M134").Value = 1447.2

$ws.Range("H140").Value = 4476.3477
$ws.Range("I140").Value = 1830.2142
$ws.Range("J140").Value = 8592.556
$ws.Range("K140").Value = 5490.642599999999
$ws.Range("L140").Value = 25777.668
$ws.Range("M140").Value = -310.6425999999992
$ws.Range("N140").Value = -36137.66800000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1842.1
$ws.Range("J107").Value = 3224.2222
$ws.Range("L107").Value = 3224.2222
$ws.Range("N107").Value = -7064.2222

$ws.Range("H132").Value = 4308.769
$ws.Range("I132").Value = 4259.5
$ws.Range("J132").Value = 4900
$ws.Range("K132").Value = 12778.5
$ws.Range("L132").Value = 14700
$ws.Range("M132").Value = -10248.5
$ws.Range("N132").Value = -19760

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2018.1818
$ws.Range("I46").Value = 2099.75
$ws.Range("J46").Value = 1800.6666
$ws.Range("K46").Value = 2099.75
$ws.Range("L46").Value = 1800.6666
$ws.Range("M46").Value = -1911.75
$ws.Range("N46").Value = -2176.6666

$ws.Range("H132").Value = 6512.3213
$ws.Range("I132").Value = 10208
$ws.Range("J132").Value = 3740.5625
$ws.Range("K132").Value = 30624
$ws.Range("L132").Value = 11221.6875
$ws.Range("M132").Value = -28094
$ws.Range("N132").Value = -16281.6875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 644.63635
$ws.Range("I113").Value = 361.55554
$ws.Range("J113").Value = 840.61536
$ws.Range("K113").Value = 1084.66662
$ws.Range("L113").Value = 2521.84608
$ws.Range("M113").Value = 1085.33338
$ws.Range("N113").Value = -6861.84608

$ws.Range("H126").Value = 3008.6667
$ws.Range("I126").Value = 2447
$ws.Range("J126").Value = 3458
$ws.Range("K126").Value = 7341
$ws.Range("L126").Value = 10374
$ws.Range("M126").Value = -4871
$ws.Range("N126").Value = -15314

$ws.Range("H136").Value = 1088.75
$ws.Range("I136").Value = 1088.75
$ws.Range("K136").Value = 3266.25
$ws.Range("M136").Value = -716.25
